# "Liste de commande dans Git" -> corrections:
#  1) Title run split in two ("...dans" / " Git") with the _GoBack bookmark
#     (Word's "last edit" marker) sitting at the split point.
#  2) "Init : " -> "Init" in bold + " : " unchanged, matching the other
#     bold command labels (Mkdir, Touch, Cat, ...).
#  3) The _GoBack bookmark is a document-wide singleton, so re-adding it at
#     the title automatically removes the stale one that used to sit after
#     the last "origin" near the end of the document.

$d = $word.ActiveDocument

# --- 1) Split the title and drop the _GoBack bookmark at the split point ---
$titlePara = $d.Paragraphs(1).Range
$splitAt = $titlePara.Start + "Liste de commande dans".Length
$splitPoint = $d.Range($splitAt, $splitAt)
$d.Bookmarks.Add("_GoBack", $splitPoint)

# --- 2) Bold the "Init" label, leaving " : " as a separate, non-bold run ---
$initPara = $d.Paragraphs(3).Range
$initWord = $d.Range($initPara.Start, $initPara.Start + "Init".Length)
$initWord.Bold = 1
